# Updates the cryptos list snapshot (prices + volume deltas) to match
# the latest GitHub Actions refresh. Row 34/35 (Filecoin/Stellar) and
# row 43/44 (WEMIXTOKEN/Frax) swap places as the ranking shifted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.347.15"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.732.09"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -3.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.12"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -4.36%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4242"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -10.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3596"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.94"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07458"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.120"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.48"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -4.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.072"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.164"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.735.80"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001064"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.92"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +5.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06017"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -10.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.78"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -4.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.063"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -6.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5239"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -5.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.369.72"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.32"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -5.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.414"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.09"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.360"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "149.95"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.931.17"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.276"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "126.62"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -5.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.736"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -7.80%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09050"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -6.48%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.585"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -6.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.47"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2150"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02279"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06134"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.026"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6383"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -5.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.184"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.412"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -4.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.862"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("E46").Value = "  -4.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.736"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5835"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -5.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.28"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.939"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06838"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -3.98%  "
